$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update updated odds on row 13 (Venados - Tapatio, MEXICO - LIGA DE EXPANSION MX)
$ws.Range("G13").Value = 1.75
$ws.Range("H13").Value = 3.8
$ws.Range("I13").Value = 3.95
$ws.Range("J13").Value = 2.3
$ws.Range("L13").Value = 4.25
$ws.Range("Q13").Value = 2.67
$ws.Range("R13").Value = 1.36
$ws.Range("V13").Value = 1.91
$ws.Range("X13").Value = 8.5
$ws.Range("Z13").Value = 14
$ws.Range("AB13").Value = 24
$ws.Range("AD13").Value = 7.4
$ws.Range("AG13").Value = 500
$ws.Range("AH13").Value = 12.5
$ws.Range("AI13").Value = 23
$ws.Range("AJ13").Value = 13.5
$ws.Range("AK13").Value = 60

# Remove the trailing URUGUAY - LIGA AUF URUGUAYA fixtures (rows 22-29),
# shrinking the used range down to A1:AS21
$ws.Range("A22:AS29").EntireRow.Delete() | Out-Null
